$d = $word.ActiveDocument

$replacements = @(
    @{Old = "Ins-13"; New = "Rephrase"},
    @{Old = "Ins-12"; New = "Articulate"},
    @{Old = "Ins-16"; New = "Elaborate"},
    @{Old = "Zhou-Ins"; New = "Zhou-instruction"},
    @{Old = "Ins-10"; New = "Plan"},
    @{Old = "Ins-17"; New = "Converse"},
    @{Old = "Reflection"; New = "Self-critique"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
